$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'0.40%"
$ws.Range("D3").Value = "'41.36"
$ws.Range("E3").Value = "'2.97%"
$ws.Range("D4").Value = "'5.710"
$ws.Range("E4").Value = "'-2.02%"
$ws.Range("D5").Value = "'0.08068"
$ws.Range("E5").Value = "'0.55%"
$ws.Range("D6").Value = "'2.054"
$ws.Range("E6").Value = "'5.40%"
$ws.Range("D7").Value = "'8.719"
$ws.Range("E7").Value = "'0.12%"
$ws.Range("D8").Value = "'4.517"
$ws.Range("E8").Value = "'-1.29%"
$ws.Range("E9").Value = "'-0.71%"
$ws.Range("D10").Value = "'0.9215"
$ws.Range("E10").Value = "'-2.50%"
$ws.Range("D11").Value = "'0.1253"
$ws.Range("E11").Value = "'0.23%"
$ws.Range("D12").Value = "'0.1944"
$ws.Range("E12").Value = "'-0.63%"
$ws.Range("D13").Value = "'8.291"
$ws.Range("E13").Value = "'-6.65%"
$ws.Range("D14").Value = "'0.09320"
$ws.Range("E14").Value = "'1.26%"
$ws.Range("D15").Value = "'0.03668"
$ws.Range("E15").Value = "'2.02%"
$ws.Range("D16").Value = "'0.1055"
$ws.Range("E16").Value = "'9.48%"
$ws.Range("D17").Value = "'0.001292"
$ws.Range("E17").Value = "'-0.69%"
$ws.Range("D18").Value = "'0.006166"
$ws.Range("E18").Value = "'-0.22%"
$ws.Range("E19").Value = "'0.35%"
$ws.Range("E20").Value = "'-1.23%"
$ws.Range("D21").Value = "'0.1415"
$ws.Range("E21").Value = "'0.72%"
$ws.Range("D22").Value = "'0.2652"
$ws.Range("E22").Value = "'9.61%"
$ws.Range("D23").Value = "'0.04430"
$ws.Range("E23").Value = "'0.41%"
$ws.Range("D24").Value = "'0.001262"
$ws.Range("E24").Value = "'-0.08%"
$ws.Range("D25").Value = "'0.004349"
$ws.Range("E25").Value = "'-0.07%"
$ws.Range("E26").Value = "'8.32%"
$ws.Range("D39").Value = "'0.02815"
$ws.Range("E39").Value = "'16.54%"
$ws.Range("D40").Value = "'0.05468"
$ws.Range("E40").Value = "'3.83%"
$ws.Range("D41").Value = "'0.007594"
$ws.Range("E41").Value = "'1.76%"
$ws.Range("D42").Value = "'0.009970"
$ws.Range("E42").Value = "'15.78%"
$ws.Range("D43").Value = "'0.1422"
$ws.Range("E43").Value = "'0.31%"
$ws.Range("D44").Value = "'0.002112"
$ws.Range("E44").Value = "'0.11%"
$ws.Range("E45").Value = "'8.59%"
$ws.Range("D46").Value = "'0.00006736"
$ws.Range("E46").Value = "'-2.59%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.39%"
$ws.Range("D48").Value = "'0.003169"
$ws.Range("E48").Value = "'0.28%"
$ws.Range("D49").Value = "'0.002281"
$ws.Range("E49").Value = "'59.85%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.39%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.39%"
